$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DQ_Metrics")

# The regenerated DQ-Metrics export dropped several columns that are no
# longer produced by the report:
#   - rdCase_no_py, orphaCase_no_py, tracerCase_no_py (previously L:N)
#   - missing_value_no_py, orphaMissing_no_py, implausible_codeLink_no_py,
#     outlier_no_py, ambiguous_rdCase_no_py, duplicateRdCase_no_py
#     (previously S:X)
# Delete the later block first so the column letters for the first
# block remain valid.
$ws.Range("S1:X1").EntireColumn.Delete()
$ws.Range("L1:N1").EntireColumn.Delete()

# The execution time recorded for this (re-)run changed.
$ws.Range("Q2").Value = 0.05
